$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.470.44'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.904.32'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4800'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4068'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08078'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.003'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.44'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.68%  '
$ws.Range('D12').Value = '1.906.55'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.957'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.078'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.03'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.06685'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001034'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').Value = '29.478.35'
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.545'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.79'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.90%  '
$ws.Range('D25').Value = '2.125.82'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.099'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.098'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.43'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.035'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09507'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.497'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('E34').Value = '  -2.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.543'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02253'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06079'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.177'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5893'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.912'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.94%  '
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07812'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.279'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.400'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5529'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.925'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '114.15'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.2952'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.27'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.96%  '
